$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.193424463272095
$ws.Range("B1").Value = 4.489976406097412
$ws.Range("C1").Value = 0.3004874587059021
$ws.Range("D1").Value = 0.19061279296875
$ws.Range("E1").Value = 0.2627407014369965
